$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range entirely so stale cells disappear.
$ws.Range("A1:B19").Clear()

# Touch the brand-new label/value strings first, in the exact order they need
# to land in the shared-string table (matches how the source workbook grew
# its sst when the rows below were authored).
$ws.Range("A1").Value = "key"
$ws.Range("B1").Value = "value"
$ws.Range("A5").Value = "profession"
$ws.Range("A6").Value = "beruf"
$ws.Range("A19").Value = "KaP"
$ws.Range("A20").Value = "INI"
$ws.Range("B23").Value = "8 schritt"
$ws.Range("A24").Value = "ausweichen"
$ws.Range("A26").Value = "has_grimoire"
$ws.Range("A23").Value = "geschwindigkeit"
$ws.Range("B6").Value = "Stadtwache"

# Identity block
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "Margot"
$ws.Range("A3").Value = "rasse"
$ws.Range("B3").Value = "Mensch"
$ws.Range("A4").Value = "alter"
$ws.Range("B4").Value = 31
$ws.Range("B5").Value = "Soldatin"

# Attributes block
$ws.Range("A8").Value = "MU"
$ws.Range("B8").Value = 14
$ws.Range("A9").Value = "KL"
$ws.Range("B9").Value = 14
$ws.Range("A10").Value = "IN"
$ws.Range("B10").Value = 10
$ws.Range("A11").Value = "CH"
$ws.Range("B11").Value = 10
$ws.Range("A12").Value = "FF"
$ws.Range("B12").Value = 10
$ws.Range("A13").Value = "GE"
$ws.Range("B13").Value = 12
$ws.Range("A14").Value = "KO"
$ws.Range("B14").Value = 10
$ws.Range("A15").Value = "KK"
$ws.Range("B15").Value = 12

# Points block
$ws.Range("A17").Value = "LP"
$ws.Range("B17").Value = 21
$ws.Range("A18").Value = "AsP"
$ws.Range("B18").Value = 4
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 12
$ws.Range("A21").Value = "AP"
$ws.Range("B21").Value = 0

# Movement / grimoire block
$ws.Range("B24").Value = 12
$ws.Range("B26").Value = 0

# Column widths (auto-widened by Excel for the longer new labels)
$ws.Columns("A").ColumnWidth = 23.5703125
$ws.Columns("B").ColumnWidth = 18.85546875

# Selection matches the edited workbook
$ws.Range("B10").Select() | Out-Null
